{"js": "// Remove the stray \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd pair)\n// that Word leaves behind after the last edit location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Change the author name in \"Ben changing things up!\" to \"Ilyas\", keeping\n// the rest of the sentence (\"  changing things up!\") as a separate run,\n// same as a real in-place edit of just the name would produce.\nconst results = context.document.body.search(\"Ben changing things up!\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    '<pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' +\n    '<w:r><w:t>Ilyas</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> changing things up!</w:t></w:r>' +\n    '</w:p></w:body></w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n  target.insertOoxml(ooxml, \"Replace\");\n} else {\n  // Fallback: only the name \"Ben\" is present (already partly edited).\n  const benResults = context.document.body.search(\"Ben\", { matchCase: true });\n  benResults.load(\"items\");\n  await context.sync();\n  if (benResults.items.length > 0) {\n    benResults.items[0].insertText(\"Ilyas\", \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the stray \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd pair)\n# that Word leaves behind after the last edit location.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Change the author name in \"Ben changing things up!\" to \"Ilyas\".\n$find = $d.Content.Find\n$find.Text = \"Ben\"\n$find.Replacement.Text = \"Ilyas\"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n"}
